$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.009.21'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '2.792.76'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''359.13'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").Value = '''109.85'
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("D7").Value = '''0.565'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.594'
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = '''40.12'
$ws.Range("E10").Value = '  -3.33%  '
$ws.Range("D11").Value = '''0.0857'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").Value = '''19.51'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = '''7.59'
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").Value = '3.230.62'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '2.794.30'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("E17").Value = '  +6.59%  '
$ws.Range("D18").Value = '51.937.46'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '''7.42'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '''3.14'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").Value = '''13.03'
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").Value = '0.0₃0984'
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").Value = '''273.94'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = '''70.29'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = '''10.20'
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.26'
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.146'
$ws.Range("E30").Value = '  +4.72%  '
$ws.Range("D31").Value = '''51.68'
$ws.Range("E31").Value = '  +1.97%  '
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("D33").Value = '''34.51'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("D34").Value = '''5.76'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("D35").Value = '''0.0848'
$ws.Range("E35").Value = '  +3.05%  '
$ws.Range("D36").Value = '''5.26'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").Value = '''18.23'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '''2.01'
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''122.60'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '''2.25'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").Value = '''22.31'
$ws.Range("E45").Value = '  -6.25%  '
$ws.Range("D46").Value = '2.089.47'
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").Value = '''2.23'
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").Value = '''0.932'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").Value = '''8.95'
$ws.Range("E51").Value = '  +0.89%  '
